$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 20 (YARLIS ANDRES PEREZ BABILONIA) becomes the last worker row once the
# rows below/above it disappear, so give it the closing bottom-border format
# that row 21 (the last row of the table) currently has.
$ws.Range("B21:J21").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)

# Remove the three workers that left: row 21 (DANIELA ISABEL CHEDRAHUY QUIROZ),
# row 18 (JUAN PABLO ESCALANTE BUELVAS) and row 16 (MELANIS BELLO GARRIDO).
# Delete bottom-up so the remaining row numbers don't shift under us.
$ws.Rows.Item(21).Delete()
$ws.Rows.Item(18).Delete()
$ws.Rows.Item(16).Delete()

# Refresh summary figures for the updated worker list (3 workers, 3 periods).
$ws.Range("E11").Value = 108720
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 3
